# Group all the shapes that make up the "breeding / non-breeding" figure
# (the connector arrow, the three ovals, the flag pictures and the three
# labels) into a single group shape, matching the manual "Group" operation
# performed in PowerPoint while building the cross-correlation figure.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$count = $s.Shapes.Count
$ids = @()
for ($i = 1; $i -le $count; $i++) {
    $ids += $i
}

$rng = $s.Shapes.Range($ids)
$g = $rng.Group()
$g.Name = "Group 34"

# Nudge the freshly created group back to the exact position PowerPoint
# left it at (a small manual drag after grouping) - 1010652, 1502447 EMU.
$g.Left = 79.57889943779529
$g.Top = 118.30291368582677
